# Applies the "Created graphs for optimisation comparason presentation" edit:
#   - Adds 2 rows to the "1D NEW" table (training-set style entries)
#   - Inserts a new "1D TRAINING" worksheet (with its own table) between
#     "1D NEW" and "2D"
#
# $excel / $wb resolve to the running application / active workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "1D NEW" sheet: append two rows (16 & 17) to the existing table.
# ---------------------------------------------------------------------
$wsNew = $wb.Worksheets.Item("1D NEW")
$loNew = $wsNew.ListObjects.Item(1)

$row16 = $loNew.ListRows.Add()
$r16 = $row16.Range
$r16.Cells.Item(1, 1).Value = "1dmockanderrors13"
$r16.Cells.Item(1, 2).Value = 200
$r16.Cells.Item(1, 3).Value = 200
$r16.Cells.Item(1, 4).Value = 3
$r16.Cells.Item(1, 5).Value = 0.3
$r16.Cells.Item(1, 6).Value = 60
$r16.Cells.Item(1, 7).Value = 1
$r16.Cells.Item(1, 8).Value = 0
$r16.Cells.Item(1, 9).Value = 0
$r16.Cells.Item(1, 10).Value = 0
$r16.Cells.Item(1, 11).Value = 30

# ---------------------------------------------------------------------
# 2) New "1D TRAINING" sheet, inserted after "1D NEW" (before "2D").
# ---------------------------------------------------------------------
$wsTraining = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsNew)
$wsTraining.Name = "1D TRAINING"

$wsTraining.Range("B2").Value = "training_set/"

$wsTraining.Range("B3").Value = "name"
$wsTraining.Range("C3").Value = "array length (pixels)"
$wsTraining.Range("D3").Value = "pixel pitch (um)"
$wsTraining.Range("E3").Value = "central frequency (THz)"
$wsTraining.Range("F3").Value = "FWHM (THz)"
$wsTraining.Range("G3").Value = "theta (arcminutes)"
$wsTraining.Range("H3").Value = "vibrations std (mm)"
$wsTraining.Range("I3").Value = "read noise (% of coherant peak)"
$wsTraining.Range("J3").Value = "averages"

$wsTraining.Range("B4").Value = "training_set1"
$wsTraining.Range("C4").Value = 200
$wsTraining.Range("D4").Value = 200
$wsTraining.Range("E4").Value = "0.1-10"
$wsTraining.Range("F4").Value = "0.1-10"
$wsTraining.Range("G4").Value = 60
$wsTraining.Range("H4").Value = 0
$wsTraining.Range("I4").Value = 0
$wsTraining.Range("J4").Value = 30
$wsTraining.Range("K4").Value = 900

$wsTraining.Range("B5").Value = "training_set2"
$wsTraining.Range("C5").Value = 200
$wsTraining.Range("D5").Value = 200
$wsTraining.Range("E5").Value = "0.1-10"
$wsTraining.Range("F5").Value = "0.1-10"
$wsTraining.Range("G5").Value = 60
$wsTraining.Range("H5").Value = 0
$wsTraining.Range("I5").Value = 20
$wsTraining.Range("J5").Value = 30
$wsTraining.Range("K5").Value = 900

$wsTraining.Range("B6").Value = "training_set3"
$wsTraining.Range("C6").Value = 100
$wsTraining.Range("D6").Value = 400
$wsTraining.Range("E6").Value = "0.1-3.16"
$wsTraining.Range("F6").Value = "0.1-3.16"
$wsTraining.Range("G6").Value = 60
$wsTraining.Range("H6").Value = 0
$wsTraining.Range("I6").Value = 20
$wsTraining.Range("J6").Value = 30
$wsTraining.Range("K6").Value = 49

$wsTraining.Range("K3").Value = "number of interferograms"
$wsTraining.Range("L3").Value = "Comments"

# Table for the new sheet.
$trainRange = $wsTraining.Range("B3:L16")
$loTrain = $wsTraining.ListObjects.Add(1, $trainRange, $null, 1)
$loTrain.Name = "Table432"
$loTrain.TableStyle = "TableStyleMedium6"

# Title formatting (merged, bold + centered) to match the other sheets.
$wsTraining.Range("B2:J2").Merge()
$titleRange = $wsTraining.Range("B2:J2")
$titleRange.Font.Bold = $true
$titleRange.HorizontalAlignment = -4108
$wsTraining.Range("K2").Font.Bold = $true

$wsTraining.Range("K7").Select()
$wsTraining.Application.ActiveWindow.ScrollRow = 3

# ---------------------------------------------------------------------
# 3) Second extra row (17) on the "1D NEW" table, added after the new
#    sheet's strings so the shared-string table order matches the
#    original authoring session.
# ---------------------------------------------------------------------
$row17 = $loNew.ListRows.Add()
$r17 = $row17.Range
$r17.Cells.Item(1, 1).Value = "1dmockanderrors14"
$r17.Cells.Item(1, 2).Value = 100
$r17.Cells.Item(1, 3).Value = 400
$r17.Cells.Item(1, 4).Value = 3
$r17.Cells.Item(1, 5).Value = 1
$r17.Cells.Item(1, 6).Value = 60
$r17.Cells.Item(1, 7).Value = 1
$r17.Cells.Item(1, 8).Value = 0
$r17.Cells.Item(1, 9).Value = 0
$r17.Cells.Item(1, 10).Value = 10
$r17.Cells.Item(1, 11).Value = 30

$wsNew.Range("L18").Select()

Write-Output "done"
